$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new shared string value into cell E2 (text cell referencing the
# new shared string "Région Morges-Aubonne").
$ws.Range("E2").Value = "Région Morges-Aubonne"

# Give column E (NomParoisse) an explicit custom width.
$ws.Columns("E:E").ColumnWidth = 29.64

# Move the active selection from E8 to F5.
$ws.Range("F5").Select()
